$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 412, shifting existing rows 412:469 down to 413:470
$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with its data
$ws.Range("A412").Value = 5
$ws.Range("B412").Value = 'Macroferia Regional de Talca'
$ws.Range("C412").Value = 'Maule'
$ws.Range("D412").Value = 45127
$ws.Range("E412").Value = 7
$ws.Range("F412").Value = 100112009
$ws.Range("G412").Value = 'Acelga'
$ws.Range("H412").Value = 'Sin especificar'
$ws.Range("I412").Value = 'Primera'
$ws.Range("J412").Value = 500
$ws.Range("K412").Value = 1600
$ws.Range("L412").Value = 1600
$ws.Range("M412").Value = 1600
$ws.Range("N412").Value = '$/docena de atados (4 kilos)'
$ws.Range("O412").Value = 'Región del Maule'
$ws.Range("P412").Value = 400
$ws.Range("Q412").Value = 4
$ws.Range("R412").Value = 'Hortaliza'
